# Edit script: remove w:hint="cs" from 3 empty-heading paragraphs' pPr/rPr,
# and restructure the final "staticmethod" paragraph: drop <w:rtl/> from its pPr/rPr,
# then append two new paragraphs ("Abstract class:" and the Persian abstract-class
# explanation, the latter carrying the _GoBack bookmark that used to sit on the
# "staticmethod" paragraph).

$d = $word.ActiveDocument

# --- 1) Paragraph "ارث بری:" -----------------------------------------------
$xml71 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="0025352E" w:rsidRDefault="0025352E" w:rsidP="0025352E"><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="Calibri"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Calibri" w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ارث بری:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(71).Range.InsertXML($xml71)

# --- 2) Paragraph "در ارث بری ما کلاس زیر مجموعه..." ------------------------
$xml72 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="0025352E" w:rsidRDefault="0025352E" w:rsidP="0025352E"><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="Calibri"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Calibri" w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>در ارث بری ما کلاس زیر مجموعه مان هست یک از بالاسری خودش.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(72).Range.InsertXML($xml72)

# --- 3) Paragraph "زمانی ک بخواهیم کاری را..." ------------------------------
$xml74 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00A31445" w:rsidRDefault="00A31445" w:rsidP="00A31445"><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="Calibri"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Calibri" w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>زمانی ک بخواهیم کاری را ک در برنامه انجام میدهیم به خود ارجاع داده شود برای اینکار از کلاس متد باید استفاده کنیم</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(74).Range.InsertXML($xml74)

# --- 4) Final "staticmethod" paragraph + two new trailing paragraphs -------
# InsertXML on (or exactly at the end of) the very last paragraph of the body
# cannot fully replace that paragraph, because Word never lets go of the
# document's final paragraph mark - it keeps a leftover paragraph using the
# *old* formatting. Work around this by splicing the full replacement in
# just before that paragraph starts (a position safely inside the previous
# paragraph's text), which pushes the untouched original paragraph after the
# freshly inserted content; then delete that now-duplicate original tail.
$p76End = $d.Paragraphs(76).Range.End
$insertionPoint = $d.Range($p76End - 1, $p76End - 1)
$xmlTail = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00076C8A" w:rsidRPr="003063C9" w:rsidRDefault="00076C8A" w:rsidP="00076C8A"><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="Calibri"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Calibri"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>staticmethod</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Calibri" w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> متدی هست داخل متد کلاس  رفتاری مثل یک فانکشن هست</w:t></w:r><w:r w:rsidR="001D01EC"><w:rPr><w:rFonts w:cs="Calibri" w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>. یک فانکشن مستقل هست. برای حل مسائل داخلی کلاس ها استفاده میشود. یک کاری برای کلاس انجام میدهد و همچین مستقل هست.</w:t></w:r></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="Calibri"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Calibri"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>Abstract class:</w:t></w:r></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="Calibri" w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Calibri" w:hint="cs"/><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>این کلاس ابسترکت برای اینکه یک کلاس از جنس کلاس شود استفاده میشود. یک جور برچسب بزن هست برای کلاس ها ک برای پایتون مشخص شود این شی از جنس کلاس هست. هر شی میتواند یک نوع مثل رشته ی عددی و کلاسی داشته باشد که برای کلاس ها از ابسترکت برای مشخص کردن کلاس ها استفاده میشود.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($xmlTail)

$dupStart = $d.Paragraphs($d.Paragraphs.Count - 1).Range.End - 1
$docEnd = $d.Content.End
$d.Range($dupStart, $docEnd).Delete()
